# SwaadSutra_Consolidated_2026-01-13.xlsx
# A new order came in (order #3, Ajay Dwarkunde / Pohe x1) and was recorded
# at the top of the "All Orders" log (row 2, pushing the existing two order
# rows down by one). The "Daily Summary" sheet's totals for 2026-01-13 are
# bumped accordingly (one more order, +30 revenue/pending).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")

# Push existing data rows (old row 2 "Pooja…", old row 3 "Anuradha N…") down
# by inserting a fresh row right under the header.
$ws.Rows.Item(2).Insert()

# New order's details.
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 2).Value = "2026-01-13 11:15"
$ws.Cells.Item(2, 3).Value = "Ajay Dwarkunde"
$ws.Cells.Item(2, 4).Value = "b-703"
# Phone number / collection date must stay text (leading apostrophe keeps
# Excel from reinterpreting these numeric-looking strings as a Number /
# serial date), matching the sheet's numberStoredAsText convention.
$ws.Cells.Item(2, 5).Value = "'8087172173"
$ws.Cells.Item(2, 6).Value = "Pohe x1"
$ws.Cells.Item(2, 7).Value = 30
$ws.Cells.Item(2, 8).Value = "NEW"
$ws.Cells.Item(2, 9).Value = "PENDING"
$ws.Cells.Item(2, 10).Value = "'2026-01-13"
$ws.Cells.Item(2, 11).Value = "18:50"
# Notes / Cancel Reason / Feedback are left blank for the new order, same
# as every other row.
$ws.Cells.Item(2, 12).Value = ""
$ws.Cells.Item(2, 13).Value = ""
$ws.Cells.Item(2, 14).Value = ""

# Daily Summary roll-up for 2026-01-13: one more order (2 -> 3) and +30 to
# both Revenue and Pending (the new order is unpaid, like the others).
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Cells.Item(2, 2).Value = 3
$ws2.Cells.Item(2, 5).Value = 75
$ws2.Cells.Item(2, 7).Value = 75
